$d = $word.ActiveDocument

$replacements = @(
    @{old="682×6="; new="224×4="},
    @{old="763×5="; new="269×7="},
    @{old="158×7="; new="934×3="},
    @{old="878×2="; new="131×4="},
    @{old="782×9="; new="565×6="},
    @{old="573×4="; new="480×4="},
    @{old="169×4="; new="771×2="},
    @{old="707×2="; new="227×5="},
    @{old="979×6="; new="358×7="},
    @{old="258×3="; new="515×4="},
    @{old="933×5="; new="750×8="},
    @{old="625×7="; new="249×7="},
    @{old="490×2="; new="622×9="},
    @{old="703×3="; new="866×7="},
    @{old="577×2="; new="138×7="},
    @{old="259×8="; new="110×3="},
    @{old="860×9="; new="637×7="},
    @{old="905×4="; new="297×4="},
    @{old="679×8="; new="135×3="},
    @{old="856×4="; new="483×4="},
    @{old="646×4="; new="630×6="},
    @{old="818×5="; new="838×7="},
    @{old="713×3="; new="986×2="},
    @{old="347×5="; new="305×6="},
    @{old="457×3="; new="376×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
